$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values (rows 2-19): add 206 to each existing value
$newValues = @(207,208,209,210,211,212,213,214,215,216,217,218,219,220,221,222,223,224)
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}

# Shrink row height for data rows 2-19 from 15 to 13.8
for ($row = 2; $row -le 19; $row++) {
    $ws.Rows.Item($row).RowHeight = 13.8
}

# Update the selection to C2:C19 with active cell C2
$ws.Range("C2:C19").Select()
